# Update the "Förändrad" (changed) date column (C) for rows 2-14
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
